# The commit regenerates household.xlsx via the JS XLSXConverter after
# fixing an "explicit model bug". The functional deltas versus the
# previous export are:
#   1. The "settings" sheet becomes the active/selected tab (was "choices").
#   2. The selection on the "settings" sheet moves from D3 to B7.
#   3. form_version (settings!B3) is bumped from 1 to the build-date style
#      version number 20130408.

$wb = $excel.ActiveWorkbook

# --- settings!B3 (form_version) : 1 -> 20130408 -----------------------
$settings = $wb.Worksheets("settings")
$settings.Range("B3").Value = 20130408

# --- selection on settings moves to B7, and settings becomes active ---
$settings.Activate()
$settings.Range("B7").Select()
